$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.012.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.67%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.875.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.80%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.41%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'318.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.44%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.33%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4364"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -5.31%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3761"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.83%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07495"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.35%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9373"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.48%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'21.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -5.31%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.864.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.92%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'6.755"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.25%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.453"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.26%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.06883"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.44%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.33%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'81.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.40%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000009062"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -5.23%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.02%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -5.37%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'27.995.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.82%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.130"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -4.13%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'11.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.84%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.123.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.38%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.038"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.93%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'152.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.27%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.85%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'5.623"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.09%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'113.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.69%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -7.80%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.09038"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.27%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.8151"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -5.88%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.822"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -6.02%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -5.20%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.964"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.91%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.002"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.20%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'Hedera"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.05518"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.37%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'TrustWalletToken"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'1.121"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.20%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01982"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.34%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.959"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.68%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.5285"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.33%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1704"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.91%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'7.023"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.98%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.794"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -6.25%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -2.45%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4900"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -5.87%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -5.50%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'106.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.16%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'NEARProtocol"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'1.678"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -5.79%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'RenderToken"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'1.910"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -13.96%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.9999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.44%  "
$ws.Range("E51").Style = "Normal"
